# Auto-generated edit script: updates market price / profit columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets, per scheduled
# market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 791.5714
$ws.Cells.Item(33, 9).Value = 791.5714
$ws.Cells.Item(33, 11).Value = 791.5714
$ws.Cells.Item(33, 13).Value = -562.5714

$ws.Cells.Item(80, 8).Value = 13889655
$ws.Cells.Item(80, 9).Value = 25000464
$ws.Cells.Item(80, 11).Value = 75001392
$ws.Cells.Item(80, 13).Value = -75000394

$ws.Cells.Item(83, 8).Value = 13889655
$ws.Cells.Item(83, 9).Value = 25000464
$ws.Cells.Item(83, 11).Value = 225004176
$ws.Cells.Item(83, 13).Value = -224999184

$ws.Cells.Item(105, 8).Value = 29000
$ws.Cells.Item(105, 10).Value = 29000
$ws.Cells.Item(105, 12).Value = 29000
$ws.Cells.Item(105, 14).Value = -35988

$ws.Cells.Item(113, 8).Value = 61497.25
$ws.Cells.Item(113, 9).Value = 3571.4285
$ws.Cells.Item(113, 10).Value = 142593.4
$ws.Cells.Item(113, 11).Value = 3571.4285
$ws.Cells.Item(113, 12).Value = 142593.4
$ws.Cells.Item(113, 13).Value = -317.4285
$ws.Cells.Item(113, 14).Value = -149101.4

$ws.Cells.Item(132, 8).Value = 3044.2856
$ws.Cells.Item(132, 9).Value = 1211.85
$ws.Cells.Item(132, 10).Value = 11188.444
$ws.Cells.Item(132, 11).Value = 3635.55
$ws.Cells.Item(132, 12).Value = 33565.33199999999
$ws.Cells.Item(132, 13).Value = -1105.55
$ws.Cells.Item(132, 14).Value = -38625.33199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 970353.25
$ws.Cells.Item(2, 9).Value = 1027403.44
$ws.Cells.Item(2, 11).Value = 1027403.44
$ws.Cells.Item(2, 13).Value = -1027290.44

$ws.Cells.Item(45, 8).Value = 3308.5
$ws.Cells.Item(45, 9).Value = 1998
$ws.Cells.Item(45, 10).Value = 4094.8
$ws.Cells.Item(45, 11).Value = 1998
$ws.Cells.Item(45, 12).Value = 4094.8
$ws.Cells.Item(45, 13).Value = -1621
$ws.Cells.Item(45, 14).Value = -4848.8

$ws.Cells.Item(46, 8).Value = 21598.375
$ws.Cells.Item(46, 9).Value = 8900
$ws.Cells.Item(46, 10).Value = 25831.166
$ws.Cells.Item(46, 11).Value = 8900
$ws.Cells.Item(46, 12).Value = 25831.166
$ws.Cells.Item(46, 13).Value = -8581
$ws.Cells.Item(46, 14).Value = -26469.166

$ws.Cells.Item(61, 8).Value = 3911.853
$ws.Cells.Item(61, 10).Value = 6040.5
$ws.Cells.Item(61, 12).Value = 6040.5
$ws.Cells.Item(61, 14).Value = -6464.5

$ws.Cells.Item(74, 8).Value = 57693456
$ws.Cells.Item(74, 9).Value = 62500870
$ws.Cells.Item(74, 11).Value = 62500870
$ws.Cells.Item(74, 13).Value = -62499996

$ws.Cells.Item(77, 8).Value = 57693456
$ws.Cells.Item(77, 9).Value = 62500870
$ws.Cells.Item(77, 11).Value = 312504350
$ws.Cells.Item(77, 13).Value = -312499982

$ws.Cells.Item(102, 8).Value = 1056012.6
$ws.Cells.Item(102, 9).Value = 1247721.9
$ws.Cells.Item(102, 10).Value = 1612
$ws.Cells.Item(102, 11).Value = 1247721.9
$ws.Cells.Item(102, 12).Value = 1612
$ws.Cells.Item(102, 13).Value = -1246099.9
$ws.Cells.Item(102, 14).Value = -4856

$ws.Cells.Item(116, 8).Value = 970353.25
$ws.Cells.Item(116, 9).Value = 1027403.44
$ws.Cells.Item(116, 11).Value = 1027403.44
$ws.Cells.Item(116, 13).Value = -1025109.44

$ws.Cells.Item(122, 8).Value = 3790447.2
$ws.Cells.Item(122, 9).Value = 4169092.2
$ws.Cells.Item(122, 11).Value = 12507276.6
$ws.Cells.Item(122, 13).Value = -12504826.6

$ws.Cells.Item(132, 8).Value = 17205.363
$ws.Cells.Item(132, 9).Value = 21315.072
$ws.Cells.Item(132, 11).Value = 63945.216
$ws.Cells.Item(132, 13).Value = -61415.216

$ws.Cells.Item(136, 8).Value = 3911.853
$ws.Cells.Item(136, 10).Value = 6040.5
$ws.Cells.Item(136, 12).Value = 18121.5
$ws.Cells.Item(136, 14).Value = -23221.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 970353.25
$ws.Cells.Item(3, 9).Value = 1027403.44
$ws.Cells.Item(3, 11).Value = 1027403.44
$ws.Cells.Item(3, 13).Value = -1027289.44

$ws.Cells.Item(86, 8).Value = 90912480
$ws.Cells.Item(86, 9).Value = 2797.1428
$ws.Cells.Item(86, 10).Value = 250004430
$ws.Cells.Item(86, 11).Value = 2797.1428
$ws.Cells.Item(86, 12).Value = 250004430
$ws.Cells.Item(86, 13).Value = -1674.1428
$ws.Cells.Item(86, 14).Value = -250006676

$ws.Cells.Item(89, 8).Value = 90912480
$ws.Cells.Item(89, 9).Value = 2797.1428
$ws.Cells.Item(89, 10).Value = 250004430
$ws.Cells.Item(89, 11).Value = 13985.714
$ws.Cells.Item(89, 12).Value = 1250022150
$ws.Cells.Item(89, 13).Value = -8369.714
$ws.Cells.Item(89, 14).Value = -1250033382

$ws.Cells.Item(105, 8).Value = 4360.2144
$ws.Cells.Item(105, 9).Value = 4420.25
$ws.Cells.Item(105, 11).Value = 4420.25
$ws.Cells.Item(105, 13).Value = -2673.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 2759.9688
$ws.Cells.Item(7, 9).Value = 2481.6086
$ws.Cells.Item(7, 10).Value = 3471.3333
$ws.Cells.Item(7, 11).Value = 2481.6086
$ws.Cells.Item(7, 12).Value = 3471.3333
$ws.Cells.Item(7, 13).Value = -2368.6086
$ws.Cells.Item(7, 14).Value = -3697.3333

$ws.Cells.Item(31, 8).Value = 12988807
$ws.Cells.Item(31, 9).Value = 14494084
$ws.Cells.Item(31, 11).Value = 14494084
$ws.Cells.Item(31, 13).Value = -14493789

$ws.Cells.Item(34, 8).Value = 12988807
$ws.Cells.Item(34, 9).Value = 14494084
$ws.Cells.Item(34, 11).Value = 14494084
$ws.Cells.Item(34, 13).Value = -14493882

$ws.Cells.Item(107, 8).Value = 910116.8
$ws.Cells.Item(107, 9).Value = 1818753.5
$ws.Cells.Item(107, 11).Value = 1818753.5
$ws.Cells.Item(107, 13).Value = -1816833.5

$ws.Cells.Item(138, 8).Value = 99812.586
$ws.Cells.Item(138, 10).Value = 99775.2
$ws.Cells.Item(138, 12).Value = 99775.2
$ws.Cells.Item(138, 14).Value = -110055.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 1127.7142
$ws.Cells.Item(92, 9).Value = 898.75
$ws.Cells.Item(92, 11).Value = 2696.25
$ws.Cells.Item(92, 13).Value = -1448.25

$ws.Cells.Item(113, 8).Value = 1251.5
$ws.Cells.Item(113, 10).Value = 1504
$ws.Cells.Item(113, 12).Value = 4512
$ws.Cells.Item(113, 14).Value = -8852

$ws.Cells.Item(131, 8).Value = 20098570
$ws.Cells.Item(131, 9).Value = 11113039
$ws.Cells.Item(131, 10).Value = 24446408
$ws.Cells.Item(131, 11).Value = 33339117
$ws.Cells.Item(131, 12).Value = 73339224
$ws.Cells.Item(131, 13).Value = -33334077
$ws.Cells.Item(131, 14).Value = -73349304

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 1521804.1
$ws.Cells.Item(70, 9).Value = 2680120.2
$ws.Cells.Item(70, 10).Value = 7082.846
$ws.Cells.Item(70, 11).Value = 2680120.2
$ws.Cells.Item(70, 12).Value = 7082.846
$ws.Cells.Item(70, 13).Value = -2679850.2
$ws.Cells.Item(70, 14).Value = -7622.846

$ws.Cells.Item(73, 8).Value = 1521804.1
$ws.Cells.Item(73, 9).Value = 2680120.2
$ws.Cells.Item(73, 10).Value = 7082.846
$ws.Cells.Item(73, 11).Value = 2680120.2
$ws.Cells.Item(73, 12).Value = 7082.846
$ws.Cells.Item(73, 13).Value = -2679184.2
$ws.Cells.Item(73, 14).Value = -8954.846

$ws.Cells.Item(122, 8).Value = 349274.12
$ws.Cells.Item(122, 9).Value = 716253.7
$ws.Cells.Item(122, 11).Value = 2148761.1
$ws.Cells.Item(122, 13).Value = -2146311.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1930
$ws.Cells.Item(100, 9).Value = 1930
$ws.Cells.Item(100, 11).Value = 1930
$ws.Cells.Item(100, 13).Value = -1389

$ws.Cells.Item(104, 8).Value = 43794.375
$ws.Cells.Item(104, 10).Value = 43794.375
$ws.Cells.Item(104, 12).Value = 43794.375
$ws.Cells.Item(104, 14).Value = -50782.375

$ws.Cells.Item(106, 8).Value = 18294
$ws.Cells.Item(106, 10).Value = 18294
$ws.Cells.Item(106, 12).Value = 18294
$ws.Cells.Item(106, 14).Value = -20818

$ws.Cells.Item(132, 8).Value = 3156.7805
$ws.Cells.Item(132, 9).Value = 2725.4026
$ws.Cells.Item(132, 10).Value = 9800
$ws.Cells.Item(132, 11).Value = 8176.2078
$ws.Cells.Item(132, 12).Value = 29400
$ws.Cells.Item(132, 13).Value = -5646.2078
$ws.Cells.Item(132, 14).Value = -34460

$ws.Cells.Item(136, 8).Value = 3943.7637
$ws.Cells.Item(136, 9).Value = 2982.1191
$ws.Cells.Item(136, 11).Value = 8946.3573
$ws.Cells.Item(136, 13).Value = -6396.3573

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2366.1853
$ws.Cells.Item(132, 9).Value = 1908.4783
$ws.Cells.Item(132, 11).Value = 5725.4349
$ws.Cells.Item(132, 13).Value = -3195.4349

$ws.Cells.Item(136, 8).Value = 6550.7144
$ws.Cells.Item(136, 9).Value = 4771
$ws.Cells.Item(136, 10).Value = 11000
$ws.Cells.Item(136, 11).Value = 14313
$ws.Cells.Item(136, 12).Value = 33000
$ws.Cells.Item(136, 13).Value = -11763
$ws.Cells.Item(136, 14).Value = -38100
